# Add the 2020 data column (X) to the "Prevalence of undernourishment" table.
# Column W (2019) is the most recent existing year column; the new column X
# extends the same header/data/formatting pattern for 2020.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting of the 2019 column (W4:W16) onto the new 2020 column
# (X4:X16) so the new cells render like the rest of the year columns.
$ws.Range("W4:W16").Copy() | Out-Null
$ws.Range("X4:X16").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

# Header (row 4) + the 2020 values for each row of the table.
$ws.Range("X4").Value = 2020
$ws.Range("X5").Value = 45.3
$ws.Range("X6").Value = 48.2
$ws.Range("X7").Value = 43.6
$ws.Range("X8").Value = 48.8
$ws.Range("X9").Value = 41.5
$ws.Range("X10").Value = 49.7
$ws.Range("X11").Value = 46.7
$ws.Range("X12").Value = 36.5
$ws.Range("X13").Value = 29.6
$ws.Range("X14").Value = 54.7
$ws.Range("X15").Value = 51.6
$ws.Range("X16").Value = 47.2

# Match the author's final cursor position/selection.
$ws.Range("AI21").Select() | Out-Null
